$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new row 7: USERNAME / PASSWORD / 1
$ws.Range("A7").Value = "USERNAME"
$ws.Range("B7").Value = "PASSWORD"
$ws.Range("C7").Value = 1

# Match the formatting of the other data rows (A2:C6 use style "1" -> centered alignment)
$ws.Range("A6:C6").Copy()
$ws.Range("A7:C7").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Update the active selection to mirror the author's saved state
$ws.Range("F7").Select()
